# Update the "YOY Expense & Profitability Analysis" table: the F8 actual
# value (Marketing & Advertising - this period) changes from 66195 to 64665.
# All of G8, F9, G9, F10, G10 are formulas that depend (directly or via
# chained formulas) on F8, so Excel recalculates them automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F8").Value = 64665

# Reflect the author's final selection/scroll state: the active cell moves
# to F10 and the sheet view no longer pins a frozen/scrolled top-left cell.
$null = $ws.Range("F10").Select()
